$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("__data")

# Row 6: Grandma Covenant -> Pocket Chronometer
$ws.Range("E6").Value = "Pocket Chronometer"
$ws.Range("F6").Value = "offlineMultiplier"
$ws.Range("G6").Value = 2
$ws.Range("H6").Value = 35
$ws.Range("I6").Value = "离线收益翻倍，归来即可收割。"
$ws.Range("J6").Value = "https://cdn.jsdelivr.net/gh/twitter/twemoji@14.0.2/assets/72x72/23f1.png"

# Row 7: Factory Time Dilation -> Quantum Ledger
$ws.Range("E7").Value = "Quantum Ledger"
$ws.Range("F7").Value = "costReduction"
$ws.Range("G7").Value = 0.15
$ws.Range("H7").Value = 40
$ws.Range("I7").Value = "所有建筑成本降低 15%。"
$ws.Range("J7").Value = "https://cdn.jsdelivr.net/gh/twitter/twemoji@14.0.2/assets/72x72/1f4b0.png"

# Row 8: Chrono Crumbs -> Starfarer Compass
$ws.Range("E8").Value = "Starfarer Compass"
$ws.Range("F8").Value = "prestigeBonus"
$ws.Range("G8").Value = 0.25
$ws.Range("H8").Value = 45
$ws.Range("I8").Value = "声望重置额外 +25% 神器点。"
$ws.Range("J8").Value = "https://cdn.jsdelivr.net/gh/twitter/twemoji@14.0.2/assets/72x72/1f320.png"
